$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text
$textCells = @("D5", "D6", "D7", "D8", "D12", "D13", "D15", "D19", "D21", "D22", "D24", "D25", "D28", "D30", "D31", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D48", "D49", "D50", "D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '60.548.93'
$ws.Range("E2").Value = '  +1.75%  '
$ws.Range("D3").Value = '2.600.72'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '583.57'
$ws.Range("E5").Value = '  +5.74%  '
$ws.Range("D6").Value = '143.38'
$ws.Range("E6").Value = '  +1.65%  '
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("D8").Value = '0.599'
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").Value = '2.622.54'
$ws.Range("E9").Value = '  +1.67%  '
$ws.Range("E10").Value = '  -2.40%  '
$ws.Range("E11").Value = '  +1.80%  '
$ws.Range("D12").Value = '0.156'
$ws.Range("E12").Value = '  -4.03%  '
$ws.Range("D13").Value = '0.375'
$ws.Range("E13").Value = '  +5.97%  '
$ws.Range("D14").Value = '3.068.22'
$ws.Range("E14").Value = '  +1.36%  '
$ws.Range("D15").Value = '24.78'
$ws.Range("E15").Value = '  +7.47%  '
$ws.Range("D16").Value = '60.555.39'
$ws.Range("E16").Value = '  +1.98%  '
$ws.Range("E17").Value = '  +3.32%  '
$ws.Range("D18").Value = '2.616.19'
$ws.Range("E18").Value = '  +1.64%  '
$ws.Range("D19").Value = '11.40'
$ws.Range("E19").Value = '  +10.59%  '
$ws.Range("E20").Value = '  +3.46%  '
$ws.Range("D21").Value = '348.68'
$ws.Range("E21").Value = '  +3.09%  '
$ws.Range("D22").Value = '6.91'
$ws.Range("E22").Value = '  +7.21%  '
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").Value = '0.522'
$ws.Range("E24").Value = '  +8.90%  '
$ws.Range("D25").Value = '63.15'
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("E26").Value = '  -0.13%  '
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '7.93'
$ws.Range("E28").Value = '  +7.07%  '
$ws.Range("E29").Value = '  +3.58%  '
$ws.Range("D30").Value = '1.86'
$ws.Range("E30").Value = '  +11.23%  '
$ws.Range("D31").Value = '6.38'
$ws.Range("E31").Value = '  +3.15%  '
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").Value = '164.74'
$ws.Range("E33").Value = '  +3.69%  '
$ws.Range("D34").Value = '19.51'
$ws.Range("E34").Value = '  +2.48%  '
$ws.Range("D35").Value = '1.01'
$ws.Range("E35").Value = '  +12.27%  '
$ws.Range("E36").Value = '  +4.08%  '
$ws.Range("D37").Value = '1.24'
$ws.Range("E37").Value = '  +5.71%  '
$ws.Range("E38").Value = '  +10.28%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '316.05'
$ws.Range("E39").Value = '  +10.17%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '37.95'
$ws.Range("E40").Value = '  +1.16%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '3.90'
$ws.Range("E41").Value = '  +6.18%  '
$ws.Range("D42").Value = '0.846'
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("D43").Value = '135.26'
$ws.Range("E43").Value = '  -0.66%  '
$ws.Range("D44").Value = '0.0992'
$ws.Range("E44").Value = '  +2.32%  '
$ws.Range("D45").Value = '0.995'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("E47").Value = '  +11.76%  '
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").Value = '0.0553'
$ws.Range("E48").Value = '  +4.37%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '0.607'
$ws.Range("E49").Value = '  +2.71%  '
$ws.Range("D50").Value = '20.23'
$ws.Range("E50").Value = '  +8.17%  '
$ws.Range("D51").Value = '0.0243'
$ws.Range("E51").Value = '  +4.20%  '
